# Petty cash book update - 25-Jan-2021, midday update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: updated opening/carried balance ---
$ws.Range("E2").Value2 = 715525

# --- Row 3: date moves forward, Wages Expense (D3) reduced ---
$ws.Range("A3").Value2 = 44221
$ws.Range("D3").Formula = "=45000"

# --- Row 4: was A/R (C4 credit), now TRANSFER BCA (D4 debit) ---
$ws.Range("B4").Value2 = "TRANSFER BCA"
$ws.Range("C4").Clear()
$ws.Range("D4").Formula = "=49000+37256000+3465000+6240000"

# --- Row 5: was TRANSFER BCA (D5 debit), now A/R (C5 credit) ---
$ws.Range("B5").Value2 = "A/R"
$ws.Range("D5").Clear()
$ws.Range("C5").Formula = "=37256000"

# --- Rows 6-34: all subsequent petty-cash entries for the period are removed ---
$ws.Range("A6:D34").Clear()

# --- Recalculate everything (the shared E-column formulas cascade the new balance) ---
$wb.Application.Calculate()

# --- View state: select D5 (frozen pane's top-left also resets to A3 on save) ---
$ws.Range("D5").Select()
